$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$beden = "32-33-34-36-38-40 Beden seçeneği mevcuttur.Ürünümüz serili olarak satılmaktadır.Belirtilen fiyatlar adet fiyatıdır."

# Row 125 - filled in D, A, B, C, E, F order (image first, then name, price, category, desc, stock)
$ws.Cells.Item(125, 4).Value = "246.jpg"
$ws.Cells.Item(125, 1).Value = "STRAİGHT FİT 246"
$ws.Cells.Item(125, 2).Value = "450 TL"
$ws.Cells.Item(125, 3).Value = "Jeans"
$ws.Cells.Item(125, 5).Value = $beden
$ws.Cells.Item(125, 6).Value = "VAR"

# Remaining rows filled in natural left-to-right order: A, B, C, D, E, F
$rows = @(
    @(126, "STRAİGHT FİT 246/3", "246-3.jpg"),
    @(127, "STRAİGHT FİT 246/4", "246-4.jpg"),
    @(128, "STRAİGHT FİT 243/1", "243-1.jpg"),
    @(129, "STRAİGHT FİT 246/2", "246-2.jpg"),
    @(130, "STRAİGHT FİT 243/4", "243-4.jpg"),
    @(131, "STRAİGHT FİT 243/5", "243-5.jpg"),
    @(132, "STRAİGHT FİT 243/2", "243-2.jpg")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = "450 TL"
    $ws.Cells.Item($r, 3).Value = "Jeans"
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $beden
    $ws.Cells.Item($r, 6).Value = "VAR"
}

$ws.Range("D140").Select()
$excel.ActiveWindow.ScrollRow = 112
